# Canteen reports update — apply per commit "Updated modules acc to LicenseKey"
$wb = $excel.ActiveWorkbook

function Set-ColWidth($sheet, $colIndex, $target) {
    $sheet.Columns.Item($colIndex).ColumnWidth = $target - (5.0/6.0)
}

# ============================================================
# Sheet 1: "Daily Summary"
# ============================================================
$ws1 = $wb.Worksheets.Item("Daily Summary")

# Collapse the 6 data rows down to 2 data rows (remove rows 4,5,6 - three middle days)
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(4).Delete()

# Column widths
Set-ColWidth $ws1 2 8
Set-ColWidth $ws1 3 7
Set-ColWidth $ws1 5 8
Set-ColWidth $ws1 6 5

# Header renames
$ws1.Cells.Item(1,2).Value = "DINNER"
$ws1.Cells.Item(1,3).Value = "LUNCH"
$ws1.Cells.Item(1,4).Value = "PIZZA"
$ws1.Cells.Item(1,5).Value = "SNACKS"
$ws1.Cells.Item(1,6).Value = "TEA"

# Row 2 (29-05-2025)
$ws1.Cells.Item(2,1).Value = "29-05-2025"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 0
$ws1.Cells.Item(2,4).Value = 0
$ws1.Cells.Item(2,5).Value = 6
$ws1.Cells.Item(2,6).Value = 4
$ws1.Cells.Item(2,7).Value = 10

# Row 3 (31-05-2025)
$ws1.Cells.Item(3,1).Value = "31-05-2025"
$ws1.Cells.Item(3,2).Value = 1
$ws1.Cells.Item(3,3).Value = 19
$ws1.Cells.Item(3,4).Value = 1
$ws1.Cells.Item(3,5).Value = 10
$ws1.Cells.Item(3,6).Value = 8
$ws1.Cells.Item(3,7).Value = 39

# Row 4 (Total)
$ws1.Cells.Item(4,2).Value = 1
$ws1.Cells.Item(4,3).Value = 19
$ws1.Cells.Item(4,4).Value = 1
$ws1.Cells.Item(4,5).Value = 16
$ws1.Cells.Item(4,6).Value = 12
$ws1.Cells.Item(4,7).Value = 49

# ============================================================
# Sheet 2: "User Summary"
# ============================================================
$ws2 = $wb.Worksheets.Item("User Summary")

# Insert a new row at position 3 (between the existing data row and the Total row),
# copying format from row 2 so styles line up.
$ws2.Rows.Item(3).Insert()
$ws2.Range("A2:I2").Copy($ws2.Range("A3:I3"))

# Column widths
Set-ColWidth $ws2 5 8
Set-ColWidth $ws2 6 7
Set-ColWidth $ws2 8 8
Set-ColWidth $ws2 9 5

# Header renames
$ws2.Cells.Item(1,5).Value = "DINNER"
$ws2.Cells.Item(1,6).Value = "LUNCH"
$ws2.Cells.Item(1,7).Value = "PIZZA"
$ws2.Cells.Item(1,8).Value = "SNACKS"
$ws2.Cells.Item(1,9).Value = "TEA"

# Row 2 (PunchID 9)
$ws2.Cells.Item(2,3).Value = 9
$ws2.Cells.Item(2,5).Value = 0
$ws2.Cells.Item(2,6).Value = 1
$ws2.Cells.Item(2,7).Value = 0
$ws2.Cells.Item(2,8).Value = 0
$ws2.Cells.Item(2,9).Value = 0

# Row 3 (PunchID 22, new row)
$ws2.Cells.Item(3,1).Value = 2025
$ws2.Cells.Item(3,2).Value = 5
$ws2.Cells.Item(3,3).Value = 22
$ws2.Cells.Item(3,5).Value = 1
$ws2.Cells.Item(3,6).Value = 18
$ws2.Cells.Item(3,7).Value = 1
$ws2.Cells.Item(3,8).Value = 16
$ws2.Cells.Item(3,9).Value = 12

# Row 4 (Total, shifted down from row 3)
$ws2.Cells.Item(4,5).Value = 1
$ws2.Cells.Item(4,6).Value = 19
$ws2.Cells.Item(4,7).Value = 1
$ws2.Cells.Item(4,8).Value = 16
$ws2.Cells.Item(4,9).Value = 12

# ============================================================
# Sheet 3: "Consumption Detail"
# ============================================================
$ws3 = $wb.Worksheets.Item("Consumption Detail")

# Add a new row 7 (TEA), copying format from row 6
$ws3.Rows.Item(7).Insert()
$ws3.Range("A6:AH6").Copy($ws3.Range("A7:AH7"))

# Column widths
Set-ColWidth $ws3 2 8
Set-ColWidth $ws3 26 4

# Row 2: PunchID 22 -> 9, BREAKFAST -> LUNCH, total 210 -> 1
$ws3.Cells.Item(2,1).Value = "9"
$ws3.Cells.Item(2,2).Value = "LUNCH"
$ws3.Cells.Item(2,3).Value = 1
$ws3.Cells.Item(2,25).ClearContents()   # Y2
$ws3.Cells.Item(2,26).ClearContents()   # Z2
$ws3.Cells.Item(2,29).ClearContents()   # AC2
$ws3.Cells.Item(2,30).ClearContents()   # AD2
$ws3.Cells.Item(2,34).Value = 1         # AH2

# Row 3: (blank) -> PunchID 22, KHAANA -> DINNER, total 4 -> 1
$ws3.Cells.Item(3,1).Value = "22"
$ws3.Cells.Item(3,2).Value = "DINNER"
$ws3.Cells.Item(3,3).Value = 1
$ws3.Cells.Item(3,25).ClearContents()   # Y3
$ws3.Cells.Item(3,34).Value = 1         # AH3

# Row 4: LUNCH total 91 -> 18
$ws3.Cells.Item(4,3).Value = 18
$ws3.Cells.Item(4,26).ClearContents()   # Z4
$ws3.Cells.Item(4,31).ClearContents()   # AE4
$ws3.Cells.Item(4,34).Value = 18        # AH4

# Row 5: SNACK -> PIZZA, total 7 -> 1
$ws3.Cells.Item(5,2).Value = "PIZZA"
$ws3.Cells.Item(5,3).Value = 1
$ws3.Cells.Item(5,25).ClearContents()   # Y5
$ws3.Cells.Item(5,34).Value = 1         # AH5

# Row 6: SNACKS total 3 -> 16
$ws3.Cells.Item(6,3).Value = 16
$ws3.Cells.Item(6,30).ClearContents()   # AD6
$ws3.Cells.Item(6,32).Value = 6         # AF6
$ws3.Cells.Item(6,34).Value = 10        # AH6

# Row 7 (new): TEA, total 12
$ws3.Cells.Item(7,1).ClearContents()    # A7 stays blank
$ws3.Cells.Item(7,2).Value = "TEA"
$ws3.Cells.Item(7,3).Value = 12
$ws3.Cells.Item(7,30).ClearContents()   # AD7 (copied from AD6=3, must clear)
$ws3.Cells.Item(7,32).Value = 4         # AF7
$ws3.Cells.Item(7,34).Value = 8         # AH7
